$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source inline-string cells)

$ws.Range("D2").Value = '30.024.69'
$ws.Range("E2").Value = '  +5.31%  '

$ws.Range("D3").Value = '1.918.58'
$ws.Range("E3").Value = '  +2.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.27'
$ws.Range("E5").Value = '  +3.60%  '

$ws.Range("E6").Value = '  -0.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5275'
$ws.Range("E7").Value = '  +3.78%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4060'
$ws.Range("E8").Value = '  +3.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08457'
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("E10").Value = '  +3.00%  '

$ws.Range("E11").Value = '  +2.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.00'
$ws.Range("E12").Value = '  +7.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.350'
$ws.Range("E13").Value = '  +1.99%  '

$ws.Range("D14").Value = '1.921.20'
$ws.Range("E14").Value = '  +2.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.366'
$ws.Range("E15").Value = '  +1.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.15'
$ws.Range("E17").Value = '  +5.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001115'
$ws.Range("E18").Value = '  +0.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06744'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.21'
$ws.Range("E20").Value = '  +2.61%  '

$ws.Range("E21").Value = '  -0.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.054'
$ws.Range("E22").Value = '  +2.11%  '

$ws.Range("D23").Value = '30.037.70'

$ws.Range("E24").Value = '  +1.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.208'
$ws.Range("E25").Value = '  -1.15%  '

$ws.Range("D26").Value = '2.142.51'
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.26'
$ws.Range("E27").Value = '  -0.96%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.10'
$ws.Range("E28").Value = '  +2.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.451'
$ws.Range("E29").Value = '  +2.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.10'
$ws.Range("E30").Value = '  +2.53%  '

$ws.Range("E31").Value = '  +3.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1063'
$ws.Range("E32").Value = '  +1.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.081'
$ws.Range("E33").Value = '  +5.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.658'
$ws.Range("E34").Value = '  +1.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02519'
$ws.Range("E35").Value = '  +2.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06597'
$ws.Range("E36").Value = '  +0.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2219'
$ws.Range("E37").Value = '  +2.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.234'
$ws.Range("E38").Value = '  +3.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.994'
$ws.Range("E39").Value = '  +1.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.192'
$ws.Range("E40").Value = '  +2.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6557'
$ws.Range("E41").Value = '  +2.46%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.243'
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.60'
$ws.Range("E43").Value = '  +4.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6176'
$ws.Range("E44").Value = '  +2.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.14'
$ws.Range("E45").Value = '  +0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.755'
$ws.Range("E46").Value = '  +1.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.066'
$ws.Range("E47").Value = '  +2.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.55'
$ws.Range("E48").Value = '  +2.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.241'
$ws.Range("E49").Value = '  +2.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.172'
$ws.Range("E50").Value = '  +3.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.51'
$ws.Range("E51").Value = '  +4.08%  '
